$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (column L) mirroring the existing K column's
# formatting, then fill in the new year's data.

# Row 3 (thin bottom-border spacer row) - just needs the same blank/bordered style
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)

# Row 4 (year header row) - copy style, then set the new year value
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value2 = 2021

# Row 6 (Mammals data row)
$ws.Range("K6").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("L6").Value2 = 7.1

# Row 7 (Birds data row)
$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial(-4122)
$ws.Range("L7").Value2 = 0.5

# Row 8 (Amphibians and Reptiles data row) - no data available, use the
# existing "-" placeholder text (same as used elsewhere in the row)
$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial(-4122)
$ws.Range("L8").Value2 = "-"

# Update the active selection to match the edited workbook
$ws.Range("N5").Select()
